$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink before shifting rows so it doesn't end up
# anchored to the wrong (shifted) cell.
$ws.Range("B15").Hyperlinks.Delete()

# Rows 8 and 9 are duplicates of rows 2 and 3 ("Affenzahn Large Friend Fox"
# and "Affenzahn Small Friend Unicorn") - remove them, shifting everything
# below up by two rows.
$ws.Range("A8:E9").EntireRow.Delete()

# The "TY Kiki Cat" row (now row 9) gets a new product image, replacing the
# broken/removed placeholder gif.
$ws.Range("B9").Value = "https://m.media-amazon.com/images/I/719xT8Gej0L._AC_SL1500_.jpg"

# Re-create the hyperlink on the "Gabby's Dollhouse Purr-ific Play Room" row,
# which is now row 13 after the two-row shift.
$ws.Hyperlinks.Add($ws.Range("B13"), "https://m.media-amazon.com/images/I/91bQnHUO2SL._AC_SL1500_.jpg")
$ws.Range("B13").Style = "Hyperlink"

# Update the current selection to match the edited workbook.
$ws.Range("B9").Select()
